$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Planilha1")

# Copy formatting (incl. date number format) from the row above into the new row
$ws.Range("A16:F16").Copy()
$ws.Range("A17:F17").PasteSpecial(-4122)  # xlPasteFormats
$excel.CutCopyMode = $false

# Add new row 17 with data for 2022-04-12
$ws.Cells.Item(17, 1).Value = 44663
$ws.Cells.Item(17, 2).Value = 0
$ws.Cells.Item(17, 3).Value = 326803
$ws.Cells.Item(17, 4).Value = 6337
$ws.Cells.Item(17, 5).Value = 29
$ws.Cells.Item(17, 6).Value = 2

# Update selection to D16
$ws.Range("D16").Select()
